# Applies the diff: remove the stray empty BH69 cell and append two new
# data rows (70 and 71) to the "Results" sheet, extending the used range
# from A1:BM69 to A1:BM71.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 69: BH69 held an empty/NaN placeholder (<v/>) with no other
#     content; the edit drops the cell entirely (matches the pattern of
#     other rows, e.g. 61/63/65/67, which omit BH altogether).
$ws.Range("BH69").ClearContents()

# --- Row 70: new M3C2 stats entry for "mov.ply-ref.ply"
$ws.Range("A70").Value = "2025-09-01 17:21:35"
$ws.Range("B70").Value = "0342-0349"
$ws.Range("C70").Value = "mov.ply-ref.ply"
$ws.Range("D70").Value = 27711
$ws.Range("E70").Value = 0.04022176252969904
$ws.Range("F70").Value = 0.08044352505939807
$ws.Range("G70").Value = 116
$ws.Range("H70").Value = 0.004186063296163978
$ws.Range("I70").Value = 0.995813936703836
$ws.Range("J70").Value = 27595
$ws.Range("K70").Value = -66.313918
$ws.Range("L70").Value = 4.678322496426
$ws.Range("M70").Value = 27230
$ws.Range("N70").Value = -76.54673100000001
$ws.Range("O70").Value = 3.851327580325
$ws.Range("P70").Value = -0.088936
$ws.Range("Q70").Value = 0.067135
$ws.Range("R70").Value = -0.002403113535060699
$ws.Range("S70").Value = -0.002479
$ws.Range("T70").Value = 0.01302056646870666
$ws.Range("U70").Value = 0.01279688229623202
$ws.Range("V70").Value = 0.009201061061786555
$ws.Range("W70").Value = 0.008914873799999999
$ws.Range("X70").Value = -0.038914
$ws.Range("Y70").Value = 0.039006
$ws.Range("Z70").Value = -0.002811117554168197
$ws.Range("AA70").Value = -0.002606
$ws.Range("AB70").Value = 0.01189272589166281
$ws.Range("AC70").Value = 0.01155571491647637
$ws.Range("AD70").Value = 0.008694125045905252
$ws.Range("AE70").Value = 0.008712498900000001
$ws.Range("AF70").Value = 3
$ws.Range("AG70").Value = 0.03906169940611999
$ws.Range("AH70").Value = "rmse"
$ws.Range("AI70").Value = 27230
$ws.Range("AJ70").Value = 8999
$ws.Range("AK70").Value = 16614
$ws.Range("AL70").Value = 286
$ws.Range("AM70").Value = 79
$ws.Range("AN70").Value = 365
$ws.Range("AO70").Value = 0.02803510410958904
$ws.Range("AP70").Value = 0.03846781035790349
$ws.Range("AQ70").Value = -0.0208045
$ws.Range("AR70").Value = -0.009398
$ws.Range("AS70").Value = 0.002544
$ws.Range("AT70").Value = 0.021122
$ws.Range("AU70").Value = 0.011942
$ws.Range("AV70").Value = -0.020469
$ws.Range("AW70").Value = -0.009423249999999999
$ws.Range("AX70").Value = 0.002283
$ws.Range("AY70").Value = 0.018149
$ws.Range("AZ70").Value = 0.01170625
$ws.Range("BA70").Value = -0.002403113535060699
$ws.Range("BB70").Value = 0.01279688229623202
$ws.Range("BC70").Value = 55192845.10587456
$ws.Range("BD70").Value = 5.055060238338491
$ws.Range("BE70").Value = 0.07339386442785875
$ws.Range("BF70").Value = -0.07048613170858216
$ws.Range("BG70").Value = -0.0002238112750802268
$ws.Range("BH70").Value = -0.2616825721550063
$ws.Range("BI70").Value = 101382149.497469
$ws.Range("BJ70").Value = 0.7945330257039492
$ws.Range("BK70").Value = 3.182318231360224
$ws.Range("BL70").Value = "data\0342-0349\python_mov.ply-ref.ply_m3c2_distances.txt"
$ws.Range("BM70").Value = "data\0342-0349\python_mov.ply-ref.ply_m3c2_params.txt"

# --- Row 71: new M3C2 stats entry for "mov-ref"
$ws.Range("A71").Value = "2025-09-01 17:31:27"
$ws.Range("B71").Value = "0342-0349"
$ws.Range("C71").Value = "mov-ref"
$ws.Range("D71").Value = 709128
$ws.Range("E71").Value = 0.1245588149878983
$ws.Range("F71").Value = 0.2491176299757966
$ws.Range("G71").Value = 95
$ws.Range("H71").Value = 0.00013396735145136
$ws.Range("I71").Value = 0.9998660326485487
$ws.Range("J71").Value = 709033
$ws.Range("K71").Value = -151.9928180000001
$ws.Range("L71").Value = 151.729350141048
$ws.Range("M71").Value = 690813
$ws.Range("N71").Value = -1125.949142
$ws.Range("O71").Value = 62.170601304942
$ws.Range("P71").Value = -0.120008
$ws.Range("Q71").Value = 0.134913
$ws.Range("R71").Value = -0.0002143663524828888
$ws.Range("S71").Value = -0.002299
$ws.Range("T71").Value = 0.01462855990518499
$ws.Range("U71").Value = 0.01462698916272617
$ws.Range("V71").Value = 0.008219008269008636
$ws.Range("W71").Value = 0.006827373
$ws.Range("X71").Value = -0.043862
$ws.Range("Y71").Value = 0.043885
$ws.Range("Z71").Value = -0.00162988991521584
$ws.Range("AA71").Value = -0.002469
$ws.Range("AB71").Value = 0.009486637003424522
$ws.Range("AC71").Value = 0.009345573310344375
$ws.Range("AD71").Value = 0.006657630840762984
$ws.Range("AE71").Value = 0.0065990526
$ws.Range("AF71").Value = 3
$ws.Range("AG71").Value = 0.04388567971555496
$ws.Range("AH71").Value = "rmse"
$ws.Range("AI71").Value = 690813
$ws.Range("AJ71").Value = 246292
$ws.Range("AK71").Value = 444509
$ws.Range("AL71").Value = 15929
$ws.Range("AM71").Value = 2291
$ws.Range("AN71").Value = 18220
$ws.Range("AO71").Value = 0.05345534160263447
$ws.Range("AP71").Value = 0.04536447157847037
$ws.Range("AQ71").Value = -0.013568
$ws.Range("AR71").Value = -0.006508
$ws.Range("AS71").Value = 0.002725
$ws.Range("AT71").Value = 0.022282
$ws.Range("AU71").Value = 0.009233
$ws.Range("AV71").Value = -0.013419
$ws.Range("AW71").Value = -0.006579
$ws.Range("AX71").Value = 0.002358
$ws.Range("AY71").Value = 0.013891
$ws.Range("AZ71").Value = 0.008937
$ws.Range("BA71").Value = -0.0002143663524828888
$ws.Range("BB71").Value = 0.01462698916272617
$ws.Range("BC71").Value = 32673381941410.1
$ws.Range("BD71").Value = 1.010085389383569
$ws.Range("BE71").Value = 0.3498892426805521
$ws.Range("BF71").Value = -0.09252700291488647
$ws.Range("BG71").Value = -0.08886902402080391
$ws.Range("BH71").Value = 1.970131977073834
$ws.Range("BI71").Value = 11637763.60543192
$ws.Range("BJ71").Value = 2.908214737699866
$ws.Range("BK71").Value = 16.58429873827257
$ws.Range("BL71").Value = "data\0342-0349\python_mov-ref_m3c2_distances.txt"
$ws.Range("BM71").Value = "data\0342-0349\python_mov-ref_m3c2_params.txt"

Write-Host ("Used range now: " + $ws.UsedRange.Address())
